# Remove "and IT experts" from the five "Responses of end-users ... on PIECES
# Software Evaluation ..." captions throughout the questionnaire document.
# Each caption lives in its own table cell / run, so anchoring on the full
# original sentence keeps every replacement unique and unambiguous.

$d = $word.ActiveDocument

$replacements = @(
    @{
        Old = "Responses of end-users and IT experts on PIECES Software Evaluation of Performance throughout and response time."
        New = "Responses of end-users on PIECES Software Evaluation of Performance throughout and response time."
    },
    @{
        Old = "B. Responses of end-users and IT experts on PIECES Software Evaluation in terms of Information Input, Output"
        New = "B. Responses of end-users on PIECES Software Evaluation in terms of Information Input, Output"
    },
    @{
        Old = ". Responses of end-users and IT experts on PIECES Software Evaluation in terms of Control and security"
        New = ". Responses of end-users on PIECES Software Evaluation in terms of Control and security"
    },
    @{
        Old = ". Responses of end-users and IT experts on PIECES Software Evaluation in terms of Efficiency"
        New = ". Responses of end-users on PIECES Software Evaluation in terms of Efficiency"
    },
    @{
        Old = ". Responses of end-users and IT experts on PIECES Software Evaluation in terms of Service"
        New = ". Responses of end-users on PIECES Software Evaluation in terms of Service"
    }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}

# The DefaultParagraphFont character style also loses its "semiHidden" flag
# in the committed revision (Word commonly drops this on a style that is
# implicitly always in use once the document is resaved). Reflect the same
# intent through the COM Style object.
$style = $d.Styles("Default Paragraph Font")
$style.Hidden = $false
